$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-8 as per regenerated save_data
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 1
